$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-09 Tuesday" "2024-01-10 Wednesday"

Replace-Text "50×49=" "29×24="
Replace-Text "96×68=" "98×94="
Replace-Text "30×98=" "49×56="
Replace-Text "61×85=" "43×55="
Replace-Text "50×12=" "80×25="
Replace-Text "36×57=" "57×58="
Replace-Text "66×90=" "29×50="
Replace-Text "74×96=" "86×32="
Replace-Text "16×70=" "54×21="
Replace-Text "39×36=" "86×45="
Replace-Text "88×54=" "13×17="
Replace-Text "50×43=" "56×28="
Replace-Text "37×38=" "29×78="
Replace-Text "76×61=" "18×84="
Replace-Text "31×96=" "70×16="
Replace-Text "73×12=" "45×99="
Replace-Text "95×89=" "18×45="
Replace-Text "28×90=" "43×78="
Replace-Text "52×74=" "39×97="
Replace-Text "97×65=" "42×56="
Replace-Text "62×31=" "15×13="
Replace-Text "64×87=" "61×46="
Replace-Text "80×42=" "21×96="
Replace-Text "34×16=" "82×19="
Replace-Text "96×69=" "91×53="
